$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59, shifting existing rows 59-60 down to 60-61
$ws.Rows.Item(59).EntireRow.Insert()

# Populate the new row 59 with the new price record
$ws.Cells.Item(59, 1).Value = 1
$ws.Cells.Item(59, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(59, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(59, 4).Value = 44858
$ws.Cells.Item(59, 5).Value = 15
$ws.Cells.Item(59, 6).Value = 100112027
$ws.Cells.Item(59, 7).Value = "Melón"
$ws.Cells.Item(59, 8).Value = "Tuna"
$ws.Cells.Item(59, 9).Value = "Segunda"
$ws.Cells.Item(59, 10).Value = 100
$ws.Cells.Item(59, 11).Value = 24000
$ws.Cells.Item(59, 12).Value = 25000
$ws.Cells.Item(59, 13).Value = 24500
$ws.Cells.Item(59, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(59, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(59, 16).Value = 1021
$ws.Cells.Item(59, 17).Value = 24
$ws.Cells.Item(59, 18).Value = "Hortaliza"
